$d = $word.ActiveDocument

# The document currently ends with a single paragraph of body text. We append
# the new "Figure legends" section (blank spacer paragraphs + the three figure
# legend paragraphs) directly as OOXML so the exact empty-paragraph / tab-stop
# / proofErr structure from the authored edit is reproduced faithfully.
$r = $d.Content
$r.Collapse(0)

$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="3765"/></w:tabs></w:pPr><w:r><w:t>Figure legends 13-Sep-23</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Fig. 1. Repertoire complexity of Masius courtship displays across the three contexts (SOLO, AUDI, COP). A) Repertoire complexity assessed as total duration (seconds; range X-Y). B) Repertoire complexity assessed as total number of elements (range X-Y). C) Repertoire complexity assessed as total number of distinct elements (range X-Y). The only clear pattern is that COP displays have a considerably smaller number of distinct elements (mostly Bow, Neck twist and ALAD).  </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Fig. 2. Syntactic complexity of Masius courtship displays across the three contexts (SOLO, AUDI, COP). A) Syntactic complexity assessed as the scaled entropy of the of the display. Low entropy (COP) connotes high predictability, of the simple COP displays. B) Syntactic complexity assessed as the compression ratio. COP displays are simple and highly compressible. C) Compression ratio vs. entropy. COP displays are compressible and have low entropy. AUDI displays have a wide range of compressibility, but generally have higher entropy than do cop displays. SOLO displays have very low compressibility and high entropy, connoting a high level of syntactic complexity of disorganized jumbles of many different display elements. The low syntactic complexity of COP displays shows that they are simple and well-organized. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Fig. 3. Masius courtship displays are more </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>similar to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> each other within contexts (SOLO, AUDI or COP) than across contexts, regardless of male identity (i.e., whether performed by the same male, or more surprisingly, when performed by different males. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>In particular, note</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> that COP displays performed by different males are more similar to each other (Fig. 3C, second bar) than are COP displays performed by the same male (Fig. 3C, first bar). Note also that this context similarity holds across each of the three contexts (panels A, B, and C). Context, therefore, is the overwhelming determinant of the syntactic organization of Masius courtship displays.</w:t></w:r></w:p>'

$r.InsertXML($newParagraphsXml)
